# Update the 取得日時 (fetched-at) timestamp on every data row of the
# "ランサーズ" sheet to the newer scrape time, leaving all other cells
# (and the "統計" sheet) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-11-11 18:33:01"
}
